{"js": "// Replace the trailing empty paragraph with two new paragraphs of text,\n// mirroring the diff:\n//   <w:p/>  ->  <w:p>Learn how to use git and github</w:p>\n//               <w:p>Now i am creating emergency branch</w:p>\nconst body = context.document.body;\n\n// The last paragraph in the body is the empty placeholder paragraph\n// that the diff turns into the first new line of text.\nconst lastParagraph = body.paragraphs.getLast();\nlastParagraph.insertText(\"Learn how to use git and github\", \"Replace\");\n\n// The second new line becomes a brand new paragraph right after it.\nlastParagraph.insertParagraph(\"Now i am creating emergency branch\", \"After\");\n\nawait context.sync();\n", "ps1": "# Replace the trailing empty paragraph with two new paragraphs of text,\n# mirroring the diff:\n#   <w:p/>  ->  <w:p>Learn how to use git and github</w:p>\n#               <w:p>Now i am creating emergency branch</w:p>\n$d = $word.ActiveDocument\n\n# The last paragraph in the body is the empty placeholder paragraph that\n# the diff turns into the first new line of text.\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.Text = \"Learn how to use git and github\"\n\n# The second new line becomes a brand new paragraph right after it.\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.InsertParagraphAfter()\n$secondParagraph = $d.Paragraphs.Last\n$secondParagraph.Range.Text = \"Now i am creating emergency branch\"\n"}
